$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$tickers = @("AAF","ABDN","ABF","ANTO","AUTO","AV","BARC","BATS","BDEV","BEZ","BF.B","BKG","BNZL","BRBY","BRK.B","BT-A","CCH","CRDA","DCC","DGE","ENT","EXPN","FCIT","FRAS","GLEN","HLMA","HSBA","HSX","IMB","INF","ITRK","JMAT","KGF","LGEN","LLOY","LSEG","AAF","ABDN","ABF","ANTO","AUTO","AV","BARC","BATS","AAF","ABDN","ABF","ANTO","AAF","ABDN","ABF","AAF","ABDN","ABF","ANTO","AUTO","AV","BARC","BATS","BDEV","BEZ","BF.B","BKG","AAF","ABDN","ABF","AAF","ABDN","ABF","ANTO","AUTO","AV","BARC","BATS","BDEV")

$startRow = 378
for ($i = 0; $i -lt $tickers.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $tickers[$i]
}

Write-Host ("Dimension: {0}" -f $ws.UsedRange.Address())
